$d = $word.ActiveDocument

# --- Title: "Play Five Pirates Online Slot for Free - Review 2021" ---
# Appears twice (the H1 heading, and the bold "meta title" run near the
# end) and both instances get the exact same new text, so a document-wide
# replace is safe and unambiguous.
$d.Content.Find.Execute("Play Five Pirates Online Slot for Free - Review 2021", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Play Five Pirates - Free Online Slot Game", 2)

# --- Meta description sentence ---
# This long sentence itself contains the substring "1,024 ways to win", so
# replace the whole sentence FIRST, before the shorter/generic
# "1,024 ways to win" replacement below, to avoid double-editing it.
$d.Content.Find.Execute("Read our Five Pirates slot game review for 2021. Play for free and enjoy 1,024 ways to win and multiple bonus features on desktop and mobile.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Read our review of Five Pirates and play this thrilling pirate-themed slot game for free.", 2)

# --- "What we like" / "What we don't like" bullet points ---
# "1,024 ways to win" also occurs verbatim inside an unrelated FAQ answer
# ("...features 5 reels and 1,024 ways to win."), which must NOT change.
# So only touch the paragraph whose entire text is exactly that bullet
# phrase, rather than doing a blind document-wide replace.
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($text -eq "1,024 ways to win") {
        $p.Range.Find.Execute("1,024 ways to win", `
            $true, $false, $false, $false, $false, $true, 1, $false, `
            "1,024 different ways to win", 2)
    }
    elseif ($text -eq "Exciting pirate atmosphere") {
        $p.Range.Find.Execute("Exciting pirate atmosphere", `
            $true, $false, $false, $false, $false, $true, 1, $false, `
            "Thrilling pirate theme", 2)
    }
    elseif ($text -eq "High-volatility potential wins") {
        $p.Range.Find.Execute("High-volatility potential wins", `
            $true, $false, $false, $false, $false, $true, 1, $false, `
            "High jackpot and potential winnings", 2)
    }
    elseif ($text -eq "No progressive jackpot") {
        $p.Range.Find.Execute("No progressive jackpot", `
            $true, $false, $false, $false, $false, $true, 1, $false, `
            "Limited free spin options", 2)
    }
}
